$p = $ppt.ActivePresentation

# Slide 8 ("Web Server") - the body text placeholder had drifted from its
# proper placeholder geometry and was overlapping the screenshot image.
# Re-set it back to the inherited/default placeholder position & size so
# the image on the slide is no longer obscured ("fixed location of image").
$s8 = $p.Slides.Item(8)
$body = $s8.Shapes.Item(2)

$body.Left = 75
$body.Top = 205
$body.Width = 874
$body.Height = 495

# Minor follow-up text touch made while fixing the layout: re-affirm the
# trailing "NodeJS" word (picked up by the spell-checker) and leave a new
# empty line at the end of the text box.
$tr = $body.TextFrame.TextRange
$nodejs = $tr.Characters(76, 6)
$nodejs.Font.Size = 36

[void]$tr.InsertAfter("`r")
